$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A dataset_id values were regenerated (rows 2-4) ---
$ws.Range("A2").Value = "nump-7qcy"
$ws.Range("A3").Value = "ebv5-4bfy"
$ws.Range("A4").Value = "cfvy-xet3"

# --- A3 no longer links out to the old dataset id; drop its hyperlink ---
# (Range.Hyperlinks.Delete clears every hyperlink on the sheet in this
# runtime, so the two that must survive - F2's and F3's attribution links -
# are re-added right after.)
$ws.Range("A3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.census.gov/programs-surveys/acs/") | Out-Null
$ws.Range("F2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.census.gov/programs-surveys/acs/") | Out-Null
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("A3").Style = "Normal"

# --- A4 (the new dataset id) now carries the hyperlink instead ---
$ws.Hyperlinks.Add($ws.Range("A4"), "https://macondointernal.demo.socrata.com/d/cfvy-xet3", $null, $null, "https://macondointernal.demo.socrata.com/d/cfvy-xet3") | Out-Null
$ws.Range("A4").Value = "cfvy-xet3"
$ws.Range("A4").Style = "Normal"

# --- View: zoomed to 160% with A4 selected ---
$excel.ActiveWindow.Zoom = 160
$ws.Range("A4").Select()
